# Update column A (row index / id) and column E (matchweek/round number)
# for rows 2..27 on the active sheet, replacing the old numeric "id" values
# in A and the "Matchweek N" text values in E with plain numeric round
# numbers, per the rodada_27 update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: row -> @(new A value, new E value)
$updates = @{
    2  = @(1, 2)
    3  = @(3, 4)
    4  = @(5, 6)
    5  = @(6, 7)
    6  = @(9, 10)
    7  = @(10, 11)
    8  = @(12, 13)
    9  = @(15, 16)
    10 = @(17, 18)
    11 = @(19, 20)
    12 = @(21, 22)
    13 = @(23, 24)
    14 = @(31, 21)
    15 = @(1, 1)
    16 = @(32, 25)
    17 = @(16, 9)
    18 = @(16, 12)
    19 = @(20, 15)
    20 = @(30, 19)
    21 = @(22, 17)
    22 = @(13, 14)
    23 = @(11, 8)
    24 = @(7, 5)
    25 = @(2, 3)
    26 = @(33, 26)
    27 = @(21, 23)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("A$row").Value = $vals[0]
    $ws.Range("E$row").Value = $vals[1]
}
